$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B) with a new timestamp
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Set "Case Sensitive" value (row 20, column B) to the text "true"
# (not the boolean TRUE). Using a formula that evaluates to the text
# "true" and then Paste-Special-Values keeps the stored type as a
# string instead of Excel auto-coercing the literal into a boolean.
$helper = $ws.Range("Z1")
$helper.Formula = "=""true"""
$helper.Copy()
$ws.Range("B20").PasteSpecial(-4163)
$helper.ClearContents()
